# Applies the rename of "aln_asym_sum_of_pairs" -> "aln_property_entropy"
# and updates the corresponding score values / conservation strings for
# rows 4 and 5 (reference_index 2 and 3), as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("H1").Value = "Metazoa_aln_property_entropy_z_score"
$ws.Range("K1").Value = "Vertebrata_aln_property_entropy_z_score"

# --- Row 4 (reference_index = 2) ---
$ws.Range("G4").Formula = '=HYPERLINK("/Users/jackson/Dropbox (MIT)/work/07-SLiM_bioinformatics/05-conservation_pipeline/examples/table_annotation/conservation_analysis/2-9606_0_002f40/2-9606_0002f40-aln_property_entropy_og_level_score_screen.png")'
$ws.Range("H4").Value = -0.9222526690491762
$ws.Range("J4").Value = "___P_______"
$ws.Range("K4").Value = -0.9222526690491762
$ws.Range("M4").Value = "___P_______"

# --- Row 5 (reference_index = 3) ---
$ws.Range("G5").Formula = '=HYPERLINK("/Users/jackson/Dropbox (MIT)/work/07-SLiM_bioinformatics/05-conservation_pipeline/examples/table_annotation/conservation_analysis/3-9606_0_002f40/3-9606_0002f40-aln_property_entropy_og_level_score_screen.png")'
$ws.Range("H5").Value = -1.112491301468793
$ws.Range("K5").Value = -1.112491301468793
